# Parts List example data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header + table data (A1:H8) -----------------------------------
$headers = @("PN","Name","Description","Supplier","Supplier PN","Pkg QTY","Pkg Price","Item")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

$rows = @(
    @("SK1001-01","Bearing","Wheel bearing","XYZ Bearing Co.","74295-942",1,2.99,"part"),
    @("SK1002-01","Board","Standard type","Skatr Dude Inc.","BRX-02",1,15.99,"part"),
    @("SK1003-01","Truck half","Truck fixed","Skatr Dude Inc.","TR1-A",1,9.87,"part"),
    @("SK1004-01","Truck half","Truck movable","Skatr Dude Inc.","TR1-B",1,12.25,"part"),
    @("SK1005-01","Truck screw","1/4-20 SHCS","Bolts R Us","92220A",50,12.86,"part"),
    @("SK1006-01","Wheel","Hard clear urethane","Skatr Dude Inc.","WHL-PRX",4,9.87,"part"),
    @("SK1007-01","Nut","1/4-20 Hex nut","Bolts R Us","95479A",50,4.88,"part")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# Carry the "left align" formatting used throughout the sheet onto the new
# table body / header cells (matches the rest of the workbook's styling).
$ws.Range("A1:H8").HorizontalAlignment = -4131

# --- Remove the now-unused trailing blank row 62 ------------------------
$ws.Rows.Item(62).Delete()

# --- Column width tweaks to fit the new data -----------------------------
$ws.Columns.Item(2).ColumnWidth = 15.28515625
$ws.Columns.Item(3).ColumnWidth = 23.140625
$ws.Columns.Item(4).ColumnWidth = 18.42578125
$ws.Columns.Item(7).ColumnWidth = 13.5703125

# --- Turn the range into an Excel Table ----------------------------------
$tableRange = $ws.Range("A1:H8")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# Record the (stale) sort-by-Name state that the table remembers.
$lo.Sort.SortFields.Clear()
$lo.Sort.SortFields.Add($ws.Range("B1:B8"))
$lo.Sort.Header = [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes

# --- Selection / cursor ---------------------------------------------------
$ws.Range("E6").Select()
